$wb = $excel.ActiveWorkbook

# --- Update conversion text on "Hoja1" ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cell = $ws1.Range("A1")
$text = $cell.Value2
$text = $text.Replace("1000 Bs = 1.84 = 6694.94 pesos", "1000 Bs = 1.85 = 6722.22 pesos")
$text = $text.Replace("6694.94 pesos = 1.83 = 947.86 Bs", "6722.22 pesos = 1.84 = 954.32 Bs")
$cell.Value = $text

# --- Update rate figures on "tasas" sheet ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 540
$ws2.Range("O10").Value = 3630
$ws2.Range("N12").Value = 3648.8
$ws2.Range("O12").Value = 518
